# DataBaru.xlsx -- "Add files via upload" re-edit.
#
# The K column (X9 / "Ketersediaan Transportasi") was re-entered: the
# K12:K28 block got reshuffled (a left-rotation by 8 rows of the original
# K12:K28 values -- K2:K11 keep their original values) and, because the
# whole K2:K28 range was re-typed/re-pasted, it lost the "0.000" number
# format it used to carry, falling back to the default General style.
# The sheet's last active selection also moved from E9 to N6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K12").Value = 2.6709999999999998
$ws.Range("K13").Value = 0.20300000000000001
$ws.Range("K14").Value = 0.255
$ws.Range("K15").Value = 3.161
$ws.Range("K16").Value = 0.70099999999999996
$ws.Range("K17").Value = 0.54800000000000004
$ws.Range("K18").Value = 3.9540000000000002
$ws.Range("K19").Value = 2.0920000000000001
$ws.Range("K20").Value = 1.002
$ws.Range("K21").Value = 0.73899999999999999
$ws.Range("K22").Value = 1.149
$ws.Range("K23").Value = 0.05
$ws.Range("K24").Value = 1.3979999999999999
$ws.Range("K25").Value = 0.35499999999999998
$ws.Range("K26").Value = 32.552
$ws.Range("K27").Value = 3.2789999999999999
$ws.Range("K28").Value = 0.48099999999999998

# The whole column (K2:K28) was re-entered without its old "0.000" number
# format, so every cell now renders with the default/General style.
$ws.Range("K2:K28").ClearFormats()

# Leave the selection where the editor left off.
$ws.Range("N6").Select()
